$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (2021-08-10 .. 2021-08-23) appended below the existing
# data (rows 2..343). Source data: date serial, nuovi pos., somma mobile
# 7gg., somma mobile 7gg. per 100mila abitanti.
$newRows = @(
    @(344, 44418, 0,  14, 109.9332548095799),
    @(345, 44419, 0,  13, 102.0808794660385),
    @(346, 44420, 5,  14, 109.9332548095799),
    @(347, 44421, 4,  17, 133.4903808402042),
    @(348, 44422, 0,  16, 125.6380054966627),
    @(349, 44423, 2,  13, 102.0808794660385),
    @(350, 44424, 2,  13, 102.0808794660385),
    @(351, 44425, 2,  15, 117.7856301531213),
    @(352, 44426, 1,  16, 125.6380054966627),
    @(353, 44427, 2,  13, 102.0808794660385),
    @(354, 44428, 1,  10, 78.52375343541422),
    @(355, 44429, 12, 22, 172.7522575579113),
    @(356, 44430, 1,  21, 164.8998822143699),
    @(357, 44431, 2,  21, 164.8998822143699)
)

# Template row (343) carries the styling used for all data rows (bold
# centered border font + date number format on column A). Copy its
# formatting down onto each new row before writing the new values.
$templateRow = $ws.Range("A343:D343")
foreach ($row in $newRows) {
    $r = $row[0]
    $dst = $ws.Range("A$r`:D$r")
    $templateRow.Copy($dst)
}

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}
